$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Source:"
$ws.Range("A12").Font.Bold = $true
$ws.Range("A13").Value = "National Corrections Reporting Program"
$ws.Range("A14").Value = "https://www.bjs.gov/index.cfm?ty=dcdetail&iid=268"

$ws.Range("A12").Select() | Out-Null
